$d = $word.ActiveDocument

# Delete the Pandoc syntax-highlighting styles, then recreate them without
# the dark "zenburn" shading (w:shd) and with the new (pygments-like) colors,
# matching pandoc --highlight-style=pygments output. Styles.Add() appends to
# the end of the styles part, and since this contiguous run is already the
# tail of word/styles.xml, deleting+recreating all of them (in original order)
# reproduces the same relative ordering.

# Delete from the end backwards: this host resolves Styles(name) to a
# positional index, and deleting earlier entries first shifts later ones,
# so we go last-to-first to keep every lookup valid.
$names = @("SourceCode", "KeywordTok", "DataTypeTok", "DecValTok", "BaseNTok", "FloatTok", "ConstantTok", "CharTok", "SpecialCharTok", "StringTok", "VerbatimStringTok", "SpecialStringTok", "ImportTok", "CommentTok", "DocumentationTok", "AnnotationTok", "CommentVarTok", "OtherTok", "FunctionTok", "VariableTok", "ControlFlowTok", "OperatorTok", "BuiltInTok", "ExtensionTok", "PreprocessorTok", "AttributeTok", "RegionMarkerTok", "InformationTok", "WarningTok", "AlertTok", "ErrorTok", "NormalTok")
for ($i = $names.Length - 1; $i -ge 0; $i--) {
    $d.Styles($names[$i]).Delete()
}

# SourceCode
$s = $d.Styles.Add("SourceCode", 1)
$s.BaseStyle = "Normal"
$s.LinkStyle = "VerbatimChar"
$s.ParagraphFormat.WordWrap = 0

# KeywordTok
$s = $d.Styles.Add("KeywordTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 2125824
$s.Font.Bold = 1

# DataTypeTok
$s = $d.Styles.Add("DataTypeTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 8336

# DecValTok
$s = $d.Styles.Add("DecValTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 7381056

# BaseNTok
$s = $d.Styles.Add("BaseNTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 7381056

# FloatTok
$s = $d.Styles.Add("FloatTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 7381056

# ConstantTok
$s = $d.Styles.Add("ConstantTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 136

# CharTok
$s = $d.Styles.Add("CharTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 10514496

# SpecialCharTok
$s = $d.Styles.Add("SpecialCharTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 10514496

# StringTok
$s = $d.Styles.Add("StringTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 10514496

# VerbatimStringTok
$s = $d.Styles.Add("VerbatimStringTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 10514496

# SpecialStringTok
$s = $d.Styles.Add("SpecialStringTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 8939195

# ImportTok
$s = $d.Styles.Add("ImportTok", 2)
$s.BaseStyle = "VerbatimChar"

# CommentTok
$s = $d.Styles.Add("CommentTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 11575392
$s.Font.Italic = 1

# DocumentationTok
$s = $d.Styles.Add("DocumentationTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 2171322
$s.Font.Italic = 1

# AnnotationTok
$s = $d.Styles.Add("AnnotationTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 11575392
$s.Font.Bold = 1
$s.Font.Italic = 1

# CommentVarTok
$s = $d.Styles.Add("CommentVarTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 11575392
$s.Font.Bold = 1
$s.Font.Italic = 1

# OtherTok
$s = $d.Styles.Add("OtherTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 2125824

# FunctionTok
$s = $d.Styles.Add("FunctionTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 8267782

# VariableTok
$s = $d.Styles.Add("VariableTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 8132377

# ControlFlowTok
$s = $d.Styles.Add("ControlFlowTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 2125824
$s.Font.Bold = 1

# OperatorTok
$s = $d.Styles.Add("OperatorTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 6710886

# BuiltInTok
$s = $d.Styles.Add("BuiltInTok", 2)
$s.BaseStyle = "VerbatimChar"

# ExtensionTok
$s = $d.Styles.Add("ExtensionTok", 2)
$s.BaseStyle = "VerbatimChar"

# PreprocessorTok
$s = $d.Styles.Add("PreprocessorTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 31420

# AttributeTok
$s = $d.Styles.Add("AttributeTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 2723965

# RegionMarkerTok
$s = $d.Styles.Add("RegionMarkerTok", 2)
$s.BaseStyle = "VerbatimChar"

# InformationTok
$s = $d.Styles.Add("InformationTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 11575392
$s.Font.Bold = 1
$s.Font.Italic = 1

# WarningTok
$s = $d.Styles.Add("WarningTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 11575392
$s.Font.Bold = 1
$s.Font.Italic = 1

# AlertTok
$s = $d.Styles.Add("AlertTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 255
$s.Font.Bold = 1

# ErrorTok
$s = $d.Styles.Add("ErrorTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 255
$s.Font.Bold = 1

# NormalTok
$s = $d.Styles.Add("NormalTok", 2)
$s.BaseStyle = "VerbatimChar"

